$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = -0.04578012498524904
$ws.Range("J3").Value = 0.5967123075441018
$ws.Range("K3").Value = 0.4408012633709308
$ws.Range("L3").Value = 2.553532143010382

$ws.Range("I20").Value = 0.1801380462221302
$ws.Range("J20").Value = 0.5767138759588911
$ws.Range("K20").Value = 0.1547621596160096
$ws.Range("L20").Value = 2.505388901929666
